$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.462.38'
$ws.Range('E2').Value = '  -1.39%  '
$ws.Range('D3').Value = '2.510.53'
$ws.Range('E3').Value = '  +7.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '297.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.588'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.82%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.548'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.73'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0802'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.38%  '
$ws.Range('D13').Value = '2.886.99'
$ws.Range('E13').Value = '  +7.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.104'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').Value = '2.527.82'
$ws.Range('E15').Value = '  +8.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.870'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.49%  '
$ws.Range('D18').Value = '45.527.65'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.62%  '
$ws.Range('D20').Value = '0.0₃0967'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '248.63'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('E24').Value = '  +1.70%  '
$ws.Range('E25').Value = '  +8.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '40.08'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +11.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.81%  '
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.78'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.76'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '148.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.48%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +21.19%  '
$ws.Range('E36').Value = '  +3.55%  '
$ws.Range('E37').Value = '  +3.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.119'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0311'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.59%  '
$ws.Range('E42').Value = '  +7.07%  '
$ws.Range('D43').Value = '2.001.81'
$ws.Range('E43').Value = '  +8.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +23.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.48%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.67%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.11%  '
$ws.Range('E50').Value = '  +4.56%  '
$ws.Range('D51').Value = '2.750.85'
$ws.Range('E51').Value = '  +8.74%  '
